$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 93
$ws.Range("K2").Value = 116
$ws.Range("B3").Value = 65
$ws.Range("E3").Value = 108
$ws.Range("J3").Value = 175
$ws.Range("B6").Value = 10
$ws.Range("B9").Value = 306
$ws.Range("C9").Value = 379
$ws.Range("D9").Value = 329
$ws.Range("E9").Value = 351
$ws.Range("F9").Value = 412
$ws.Range("H9").Value = 369
$ws.Range("I9").Value = 415
$ws.Range("J9").Value = 331
$ws.Range("L9").Value = 363
$ws.Range("B10").Value = 1038
$ws.Range("C10").Value = 1252
$ws.Range("D10").Value = 1414
$ws.Range("E10").Value = 1733
$ws.Range("F10").Value = 1748
$ws.Range("H10").Value = 456
$ws.Range("I10").Value = 699
$ws.Range("J10").Value = 577
$ws.Range("K10").Value = 566
$ws.Range("L10").Value = 539
$ws.Range("B11").Value = 1456
$ws.Range("C11").Value = 1764
$ws.Range("D11").Value = 1935
$ws.Range("E11").Value = 2260
$ws.Range("F11").Value = 2337
$ws.Range("H11").Value = 1039
$ws.Range("I11").Value = 1401
$ws.Range("J11").Value = 1200
$ws.Range("K11").Value = 1293
$ws.Range("L11").Value = 1228

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("B7").Value = 5
$ws.Range("B9").Value = 14

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("C7").Value = 27
$ws.Range("D8").Value = 39
$ws.Range("C9").Value = 75
$ws.Range("D9").Value = 82

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value = 6
$ws.Range("F8").Value = 12
$ws.Range("F10").Value = 43
$ws.Range("K10").Value = 44

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J3").Value = 9
$ws.Range("C7").Value = 26
$ws.Range("C8").Value = 49
$ws.Range("C9").Value = 80
$ws.Range("J9").Value = 45

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("D8").Value = 37
$ws.Range("J8").Value = 51
$ws.Range("C9").Value = 240
$ws.Range("D9").Value = 413
$ws.Range("E9").Value = 499
$ws.Range("I9").Value = 163
$ws.Range("J9").Value = 91
$ws.Range("C10").Value = 285
$ws.Range("D10").Value = 476
$ws.Range("E10").Value = 565
$ws.Range("I10").Value = 273
$ws.Range("J10").Value = 188

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("J5").Value = 15
$ws.Range("F6").Value = 20
$ws.Range("F7").Value = 29
$ws.Range("J7").Value = 25

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("H6").Value = 9
$ws.Range("H7").Value = 22

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("C3").Value = 3
$ws.Range("F5").Value = 29
$ws.Range("J5").Value = 25
$ws.Range("I6").Value = 6
$ws.Range("L6").Value = 12
$ws.Range("J7").Value = 12
$ws.Range("C8").Value = 79
$ws.Range("E8").Value = 88
$ws.Range("K8").Value = 57
$ws.Range("F19").Value = 43
$ws.Range("K19").Value = 44
$ws.Range("B21").Value = 14
$ws.Range("F23").Value = 17
$ws.Range("F27").Value = 24
$ws.Range("B28").Value = 80
$ws.Range("I28").Value = 76
$ws.Range("E29").Value = 21
$ws.Range("C32").Value = 75
$ws.Range("D32").Value = 82
$ws.Range("B35").Value = 16
$ws.Range("E35").Value = 19
$ws.Range("F35").Value = 13
$ws.Range("C36").Value = 80
$ws.Range("J36").Value = 45
$ws.Range("D41").Value = 18
$ws.Range("E41").Value = 19
$ws.Range("C43").Value = 12
$ws.Range("J43").Value = 6
$ws.Range("E47").Value = 62
$ws.Range("H47").Value = 34
$ws.Range("H50").Value = 22
$ws.Range("C53").Value = 285
$ws.Range("D53").Value = 476
$ws.Range("E53").Value = 565
$ws.Range("I53").Value = 273
$ws.Range("J53").Value = 188
$ws.Range("B62").Value = 24
$ws.Range("H62").Value = 11
$ws.Range("F68").Value = 36
$ws.Range("D74").Value = 69
$ws.Range("L74").Value = 15
$ws.Range("B76").Value = 41
$ws.Range("I77").Value = 70
$ws.Range("E78").Value = 42
$ws.Range("K81").Value = 11
$ws.Range("E87").Value = 28
$ws.Range("E95").Value = 75
$ws.Range("F95").Value = 54
$ws.Range("I95").Value = 20
$ws.Range("C96").Value = 14
$ws.Range("L97").Value = 5
$ws.Range("B98").Value = 12
$ws.Range("B99").Value = 1456
$ws.Range("C99").Value = 1764
$ws.Range("D99").Value = 1935
$ws.Range("E99").Value = 2260
$ws.Range("F99").Value = 2337
$ws.Range("H99").Value = 1039
$ws.Range("I99").Value = 1401
$ws.Range("J99").Value = 1200
$ws.Range("K99").Value = 1293
$ws.Range("L99").Value = 1228

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("D6").Value = 14
$ws.Range("E6").Value = 15
$ws.Range("D7").Value = 18
$ws.Range("E7").Value = 19

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K5").Value = 4
$ws.Range("K6").Value = 11

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("E8").Value = 20
$ws.Range("E9").Value = 28

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("E5").Value = 36
$ws.Range("E6").Value = 42

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("B4").Value = 1
$ws.Range("I8").Value = 38
$ws.Range("B9").Value = 80
$ws.Range("I9").Value = 76

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("E7").Value = 51
$ws.Range("H7").Value = 12
$ws.Range("E8").Value = 62
$ws.Range("H8").Value = 34

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("E3").Value = 3
$ws.Range("E9").Value = 21

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("B7").Value = 38
$ws.Range("B8").Value = 41

$ws = $wb.Worksheets.Item('River North')
$ws.Range("L5").Value = 3
$ws.Range("D6").Value = 59
$ws.Range("D7").Value = 69
$ws.Range("L7").Value = 15

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("B5").Value = 5
$ws.Range("E6").Value = 14
$ws.Range("F6").Value = 9
$ws.Range("B7").Value = 16
$ws.Range("E7").Value = 19
$ws.Range("F7").Value = 13

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L5").Value = 3

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L7").Value = 5

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("F6").Value = 16
$ws.Range("F7").Value = 24

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("H6").Value = 6
$ws.Range("B7").Value = 19
$ws.Range("B8").Value = 24
$ws.Range("H8").Value = 11

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("F7").Value = 11
$ws.Range("F8").Value = 17

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I8").Value = 26
$ws.Range("I10").Value = 70

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("E5").Value = 6
$ws.Range("I5").Value = 7
$ws.Range("F6").Value = 48
$ws.Range("E7").Value = 75
$ws.Range("F7").Value = 54
$ws.Range("I7").Value = 20

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("C6").Value = 12
$ws.Range("C7").Value = 14

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J3").Value = 1
$ws.Range("J7").Value = 12

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("B3").Value = 1

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("B7").Value = 12

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("F7").Value = 34
$ws.Range("F8").Value = 36

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("I4").Value = 2
$ws.Range("L5").Value = 4
$ws.Range("I6").Value = 6
$ws.Range("L6").Value = 12

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("C8").Value = 45
$ws.Range("E8").Value = 49
$ws.Range("K8").Value = 20
$ws.Range("C9").Value = 79
$ws.Range("E9").Value = 88
$ws.Range("K9").Value = 57

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J2").Value = 3
$ws.Range("C6").Value = 8
$ws.Range("C7").Value = 12
$ws.Range("J7").Value = 6

$ws = $wb.Worksheets.Item('Andersonville')
$ws.Range("B3").Value = 3
$ws.Range("B4").Value = 3
